$wb = $excel.ActiveWorkbook

# "Chart" sheet holds the GSC export time series in columns A:C.
$ws = $wb.Worksheets.Item("Chart")

$ws.Cells.Item(81, 1).NumberFormat = "@"
$ws.Cells.Item(81, 1).Value = "2025-12-24"
$ws.Cells.Item(81, 1).ClearFormats()
$ws.Cells.Item(81, 2).Value = 0
$ws.Cells.Item(81, 3).Value = 31
